# Revert "correccion cola de prioridad"
# This reverts commit 726e2f03b85b65654b3378c7967357ec6dd6808a.
#
# The prior commit had changed cell D13 (sheet "Hoja1") from the text
# value "21" to the text value "15". This reverts that single cell back
# to "21" (stored as text, matching D12 and the rest of the "edad"
# column, which are all shared-string text values rather than numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D13")

# Force text storage (so "21" is written as a shared string, t="s",
# rather than being auto-coerced into a numeric cell) ...
$cell.NumberFormat = "@"
$cell.Value = "21"
# ... then restore the cell's style to the default "Normal" style so the
# cell carries no explicit formatting, matching the original workbook.
$cell.Style = "Normal"
